$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "high"
$ws.Range("B14").Value = "When no event is selected (stats in default mode) and an event is edied in day-view, the stats card behaves correctly"
$ws.Range("C14").Value = "pending"

$ws.Columns.Item(2).ColumnWidth = 106.2

$ws.Range("D14").Select() | Out-Null
